$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The bill-of-quantities table is being split: the old "light / fan / exhaust
# point" rewiring filler row (row 8) is replaced by a real "Short point"
# line item, a brand new "Long point" line item is inserted right after the
# existing "Medium point" row, and all the rows below shift down by one.
# Several quantity / amount totals are recalculated to reflect the new data.
# ---------------------------------------------------------------------------

# Insert a new row at row 10 - this pushes the old row 10 ("Rewiring of 3/5
# pin ...") and everything below it down by one row, making room for the new
# "Long point" line item.
$ws.Rows("10:10").Insert()

# --- Row 8: was the empty "Rewiring of light point..." filler, now the
#     "Short point (up to 3 mtr.)" P. point line item. ---
$ws.Range("A8").Value = "P. point"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 19
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2"
$ws.Range("E8").Value = "Short point (up to 3 mtr.)"
$ws.Range("F8").Value = 256
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "4864.00"
$ws.Range("H8").Value = 0
$ws.Range("I8").NumberFormat = "@"
$ws.Range("I8").Value = ""

# --- Row 9: "Medium point (up to 6 mtr.)" stays in place, quantity/amount
#     updated. ---
$ws.Range("C9").Value = 94
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "44368.00"

# --- Row 10 (newly inserted): "Long point (up to 10 mtr.)" P. point line
#     item. ---
$ws.Range("A10").Value = "P. point"
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 92
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "4"
$ws.Range("E10").Value = "Long point  (up to 10 mtr.)"
$ws.Range("F10").Value = 662
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "60904.00"
$ws.Range("H10").Value = 0
$ws.Range("I10").NumberFormat = "@"
$ws.Range("I10").Value = ""

# --- Row 11 (was row 10): "Rewiring of 3/5 pin ..." filler row - only the
#     "Qty executed upto date" count changes. ---
$ws.Range("C11").Value = 26

# --- Row 12 (was row 11): "On board" - quantity and amount updated. ---
$ws.Range("C12").Value = 12
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "1632.00"

# --- Row 13 (was row 12): "P & F ISI marked (IS:3854) ..." switch line -
#     quantity and amount updated. ---
$ws.Range("C13").Value = 37
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "851.00"

# --- Row 14 (was row 13): "Total" row - quantity updated. ---
$ws.Range("C14").Value = 81

# --- Row 15 (was row 14): "Add Tender Premium" row - quantity updated. ---
$ws.Range("C15").Value = 56

# --- Row 16 (was row 15): "Grand Total" row - quantity updated. ---
$ws.Range("C16").Value = 17

# --- Row 18 (was row 17): "Grand Total Rs." summary - amounts updated. ---
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "112619.00"
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = "112619.00"

# --- Row 20 (was row 19): "NET PAYABLE AMOUNT Rs." summary - amounts
#     updated (Tender Premium @ 0% row, now row 19, keeps its 0.00 values). ---
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "112619.00"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "112619.00"
